$d = $word.ActiveDocument

# Helper: find the index (1-based) of the Document.Paragraphs entry that
# fully contains a given [start,end) character range.
function Get-ParaIndexForRange($startPos, $endPos) {
    $idx = -1
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $pp = $d.Paragraphs.Item($i)
        if ($pp.Range.Start -le $startPos -and $pp.Range.End -ge $endPos) {
            $idx = $i
        }
    }
    return $idx
}

# ---------------------------------------------------------------------------
# Change 1: remove the whole "If user signups with invalid phone number..."
# paragraph (ListParagraph, numId=2, red text) in its entirety - including
# its paragraph mark, so the <w:p> element disappears completely, leaving
# the two surrounding blank ListParagraph paragraphs intact.
# ---------------------------------------------------------------------------
$findPhone = $d.Content.Find
$foundPhone = $findPhone.Execute("If user signups with invalid phone number, an error message is displayed.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundPhone) {
    $fr = $findPhone.Parent
    $targetIndex = Get-ParaIndexForRange $fr.Start $fr.End
    $targetPara = $d.Paragraphs.Item($targetIndex)
    $prevPara = $d.Paragraphs.Item($targetIndex - 1)
    $killRange = $d.Range($prevPara.Range.End, $targetPara.Range.End)
    $killRange.Delete()
}

# ---------------------------------------------------------------------------
# Change 2: split the run "If user tries to signup with email address..."
# so the word "signup" is wrapped in proofErr spellStart/spellEnd markers,
# matching what Word's proofing pass would insert, without touching the
# paragraph's own formatting (pPr / numPr / rPr) or identity attributes.
# ---------------------------------------------------------------------------
$findSignup = $d.Content.Find
$foundSignup = $findSignup.Execute("If user tries to signup with email address already in use, then the user receives an error message is displayed stating the email is already in use.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundSignup) {
    $fr2 = $findSignup.Parent
    $target2 = $d.Range($fr2.Start, $fr2.End)
    $newRunsXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="45AC83E7" w14:textId="6F55DD4E" w:rsidR="00D166C4" w:rsidRPr="00A13819" w:rsidRDefault="003640E0" w:rsidP="00D166C4"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve">If user tries to </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>signup</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve"> with email address already in use, then the user receives an error message is displayed stating the email is already in use.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $target2.InsertXML($newRunsXml)
}

# ---------------------------------------------------------------------------
# Change 3: move <w:lastRenderedPageBreak/> from the "Vitals page" heading
# run down to the "If measurement is not updated..." run.
# ---------------------------------------------------------------------------
$findVitals = $d.Content.Find
$foundVitals = $findVitals.Execute("Vitals page", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundVitals) {
    $fr3 = $findVitals.Parent
    $target3 = $d.Range($fr3.Start, $fr3.End)
    $vitalsXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="079262FB" w14:textId="2E7BC18B" w:rsidR="003646B1" w:rsidRDefault="003646B1" w:rsidP="003646B1"><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/></w:rPr><w:t>Vitals page</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $target3.InsertXML($vitalsXml)
}

$findMeasurement = $d.Content.Find
$foundMeasurement = $findMeasurement.Execute("If measurement is not updated user can click on refresh button, which should update the readings.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundMeasurement) {
    $fr4 = $findMeasurement.Parent
    $target4 = $d.Range($fr4.Start, $fr4.End)
    $measurementXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="51A375FC" w14:textId="230AE7D1" w:rsidR="003646B1" w:rsidRPr="003646B1" w:rsidRDefault="003646B1" w:rsidP="003646B1"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:lastRenderedPageBreak/><w:t>If measurement is not updated user can click on refresh button, which should update the readings.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $target4.InsertXML($measurementXml)
}
